$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.242581963539124
$ws.Range("B1").Value = 2.34207010269165
$ws.Range("C1").Value = 3.390633583068848
$ws.Range("D1").Value = 2.224164247512817
$ws.Range("E1").Value = 1.373640537261963
